$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Counts" column (M) for rows 201-244 from 0 to 1
$ws.Range("M201:M244").Value = 1
